# Update attendance flags on Sheet1.
# Columns: A=Date, B=Roll, C=Name, D=Total Attendance Count, E=Real,
#          F=Duplicate, G=Invalid, H=Absent
# Set the appropriate cells from 0 to 1 for each attendance row (3-18).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: marked Invalid + Absent
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Row 4: marked present (Total Attendance Count + Real)
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1

# Row 5: marked present
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1

# Row 6: marked Absent
$ws.Range("H6").Value = 1

# Row 7: marked Absent
$ws.Range("H7").Value = 1

# Row 8: marked Absent
$ws.Range("H8").Value = 1

# Row 9: marked present
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 1

# Row 10: marked Absent
$ws.Range("H10").Value = 1

# Row 11: marked present
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 1

# Row 12: marked Absent
$ws.Range("H12").Value = 1

# Row 13: marked Absent
$ws.Range("H13").Value = 1

# Row 14: marked Absent
$ws.Range("H14").Value = 1

# Row 15: marked Absent
$ws.Range("H15").Value = 1

# Row 16: marked Absent
$ws.Range("H16").Value = 1

# Row 17: marked Absent
$ws.Range("H17").Value = 1

# Row 18: marked Absent
$ws.Range("H18").Value = 1
